# Add a "users" column to the "project hours" sheet, listing the users
# associated with each project (commit message: "add user list to project").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New header cell E1, formatted like the existing header cells (B1:D1).
$h = $ws.Cells.Item(1, 5)
$h.Value = "users"
$h.Font.Bold = $true
$h.HorizontalAlignment = -4108
$h.VerticalAlignment = -4160
$h.Borders.LineStyle = 1

# Per-project user lists, one row per existing project (rows 2-6).
$users = @(
    "['Hamidreza Jafarnejadsani', 'Gabriel Barsi Haberfeld', 'Arun Lakshmanan']",
    "['Alexander Hill', 'Alex Hill']",
    "['Karun Koppula', 'Sierra Young']",
    "['Harshal Maske']",
    "['Nicole Chan']"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $users[$i]
}
